$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-01-30 Tuesday" "2024-01-31 Wednesday"

Replace-Text "922×5=4610" "402×4=1608"
Replace-Text "470×2=940" "347×2=694"
Replace-Text "104×3=312" "800×6=4800"
Replace-Text "724×9=6516" "295×4=1180"
Replace-Text "761×7=5327" "661×4=2644"
Replace-Text "615×8=4920" "308×2=616"
Replace-Text "380×2=760" "921×5=4605"
Replace-Text "568×3=1704" "923×9=8307"
Replace-Text "176×5=880" "543×5=2715"
Replace-Text "339×5=1695" "680×5=3400"
Replace-Text "258×2=516" "914×7=6398"
Replace-Text "726×7=5082" "242×8=1936"
Replace-Text "205×7=1435" "512×2=1024"
Replace-Text "589×5=2945" "902×2=1804"
Replace-Text "959×3=2877" "618×4=2472"
Replace-Text "854×5=4270" "355×2=710"
Replace-Text "189×8=1512" "787×4=3148"
Replace-Text "645×7=4515" "438×9=3942"
Replace-Text "786×7=5502" "291×3=873"
Replace-Text "395×5=1975" "624×6=3744"
Replace-Text "426×4=1704" "714×5=3570"
Replace-Text "419×4=1676" "163×8=1304"
Replace-Text "484×8=3872" "613×3=1839"
Replace-Text "372×4=1488" "692×2=1384"
Replace-Text "487×4=1948" "445×8=3560"

Write-Output "done"
